$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9078458547592163
$ws.Range("B1").Value = 2.995015859603882
$ws.Range("C1").Value = 4.325307846069336
$ws.Range("D1").Value = 3.015062093734741
$ws.Range("E1").Value = 1.391870021820068
